# Applies the "Add files via upload" revision:
#  1. Slide 6 (sldId 263), shape id=5 "Content Placeholder 2": the last
#     bullet's text grows from "Figure shows the Boxplots after removing
#     the " to "...the outliers".
#  2. Every auto-updating "Date Placeholder" (on the Slide Master and all
#     17 Custom Layouts) gets its cached date text refreshed from
#     30-Apr-22 to 03-May-22, reflecting the later save date.

$p = $ppt.ActivePresentation

# --- 1. Slide 6 body text edit -------------------------------------------
$s = $p.Slides.Item(6)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Id -eq 5) {
        $tr = $shp.TextFrame.TextRange
        $target = $tr.Paragraphs(4, 1)
        # Shrink first so the grow-back assigns the whole new string to a
        # single run instead of diff-appending a second run.
        $target.Text = "x"
        $target2 = $shp.TextFrame.TextRange.Paragraphs(4, 1)
        $target2.Text = "Figure shows the Boxplots after removing the outliers"
    }
}

# --- 2. Refresh cached "datetimeFigureOut" placeholder text -------------
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "30-Apr-22") {
                $tr.Text = "03-May-22"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    Update-DatePlaceholder $layout.Shapes
}
